$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AP Status")

# Add new row 4 with values, matching the bold/centered/bordered style
# used by the existing data rows (same style as row 3).
$ws.Range("A4").Value = 282
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 358
$ws.Range("D4").Value = 15
$ws.Range("E4").Value = 155
$ws.Range("F4").Value = 1

# Copy the style of row 3 onto the new row 4 so it matches (s="1").
$ws.Range("A3:F3").Copy()
$ws.Range("A4:F4").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
